$d = $word.ActiveDocument

# Locate the existing author-credit paragraph ("ΛΙΟΠΕΤΑ ΔΗΜΗΤΡΑ, 1054373, 7ο ΕΤΟΣ")
# and insert a brand-new empty paragraph immediately before it.
# InsertParagraphBefore() on a collapsed Range clones the bordered paragraph
# formatting (pBdr / rPr) from the paragraph the range sits in, which is
# exactly the look-and-feel the other author-credit lines already use.
$rng = $d.Content
[void]$rng.Find.Execute("ΛΙΟΠΕΤΑ ΔΗΜΗΤΡΑ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$rng.InsertParagraphBefore()

# $rng now sits, collapsed, inside the freshly inserted (still empty) paragraph.
$rng.InsertAfter("ΔΕΜΟΣ ΔΗΜΗΤΡΗΣ,1051329 , 8ο ΕΤΟΣ")
